$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price (D) and 1h volume change (E) columns to the latest scrape.
# Price-column text is forced via NumberFormat "@" so numeric-looking strings
# (e.g. "1.000", "41.90") are preserved verbatim instead of being parsed as
# numbers; the temporary format is then cleared back to the default style.
$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "28.457.75"
$c.Style = "Normal"
$ws.Range("E2").Value = "  +0.00%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "1.824.41"
$c.Style = "Normal"
$ws.Range("E3").Value = "  -0.19%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("E5").Value = "  -0.95%  "
$ws.Range("E6").Value = "  -0.01%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.5107"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -4.26%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.3931"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -2.97%  "
$ws.Range("E9").Value = "  +0.90%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "41.90"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +0.21%  "
$ws.Range("E11").Value = "  +0.12%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "21.04"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +0.70%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "6.263"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -1.39%  "
$ws.Range("E14").Value = "  +0.03%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "7.520"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -0.30%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "1.821.53"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +0.00%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "93.10"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +4.22%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "0.00001117"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +4.08%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "0.06641"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +0.43%  "
$ws.Range("E20").Value = "  +0.89%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "1.000"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -0.06%  "
$ws.Range("E22").Value = "  +0.61%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "28.493.07"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +0.03%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "11.27"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -0.19%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "2.255"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +4.55%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "21.35"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +3.80%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "2.031.93"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +0.08%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "155.75"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -0.65%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "2.408"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -2.89%  "
$ws.Range("E30").Value = "  +0.99%  "
$ws.Range("E31").Value = "  +0.61%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "1.111"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -1.08%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "5.684"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -0.04%  "
$ws.Range("E34").Value = "  -0.16%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.07070"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -1.33%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.2216"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -2.20%  "
$ws.Range("E37").Value = "  -0.80%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "5.187"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -0.31%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.6282"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +0.11%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "11.22"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -0.84%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "1.174"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -0.67%  "
$ws.Range("E43").Value = "  -0.02%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "1.391"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -0.49%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "13.47"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +0.20%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "3.731"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +0.71%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.5898"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +0.75%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "124.20"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -1.56%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "1.989"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -0.09%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "1.195"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -0.22%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.06900"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +0.04%  "
